$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "'2027-07-29"
$ws.Range("B2").Style = "Normal"

$ws.Range("B3").Value = "'2027-07-29"
$ws.Range("B3").Style = "Normal"

$ws.Range("E5").Value = "❌ EXPIRED 3762 days ago"
$ws.Range("E6").Value = "⚠️ Expires in 17 days"
$ws.Range("E8").Value = "⚠️ Expires in 17 days"
